$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: adjust height ---
$ws.Rows.Item(4).RowHeight = 154.2

# --- Row 5: new data row ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "ReturnManagement"
$ws.Range("C5").Value = "BRD"
$ws.Range("D5").Value = "NA"

$customerFiles = "1.customer-login.coponent(Ng)`n2.customer.ts`n3.customer.service.ts`n4customer-registration component(Ng)`n5 customer product component(Ng)`n6 customer prodlist component(Ng)`n7 customer order component(Ng)`n8.Controller-webcontroller.java`n9.Service-Customer_Service.java,Customer_Sevice_Imp.java`n10-Model-Customer.java`n11.DAO-Customer_DAO.java,Customer_DAO_Imp.java`n12.Controller-webcontroller.java`n13.Service-CustomerProduct_Service.java,CustomerProduct_Sevice_Imp.java`n14-Model-CustomerProduct.java`n15.DAO-CustomerProduct_DAO.java,CustomerProduct_DAO_Imp.java`n16.Controller-webcontroller.java`n17.Service-CustomerOrder_Service.java,CustomerOrder_Sevice_Imp.java`n18-Model-CustomerOrder.java`n19.DAO-CustomerOrder_DAO.java,CustomerOrder_DAO_Imp.java`n"
$ws.Range("F5").Value = $customerFiles
$ws.Range("G5").Value = "customer ,customerproduct ,customerorder"

# Row5 specific alignments
$ws.Range("B5").WrapText = $true
$ws.Range("B5").HorizontalAlignment = -4131

$ws.Range("F5").WrapText = $true
$ws.Range("F5").VerticalAlignment = -4160

# Row height must be set after content/wrap changes so it isn't auto recalculated away
$ws.Rows.Item(5).RowHeight = 408.6

# --- View changes ---
$excel.ActiveWindow.Zoom = 64
$ws.Range("B1").Select()
